# Apply updated values to rows 2-11 of the active worksheet.
# Columns: A=Id, B=Taxonsorteringsordning, D=Rodlistade, E=TaxonId, F=Artnamn,
#          G=Vetenskapligt namn, H=Auktor, Q=Ost, R=Nord

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 112327588
$ws.Range("B2").Value = 96652
$ws.Range("Q2").Value = 641133

# Row 3
$ws.Range("A3").Value = 112327353
$ws.Range("B3").Value = 77650
$ws.Range("Q3").Value = 641135
$ws.Range("R3").Value = 7163078

# Row 4
$ws.Range("A4").Value = 112327346
$ws.Range("B4").Value = 77650
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("Q4").Value = 641252

# Row 5
$ws.Range("A5").Value = 112327585
$ws.Range("B5").Value = 96652
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 219790
$ws.Range("F5").Value = "Fläcknycklar"
$ws.Range("G5").Value = "Dactylorhiza maculata"
$ws.Range("H5").Value = "(L.) Soó"
$ws.Range("Q5").Value = 641127
$ws.Range("R5").Value = 7163079

# Row 6
$ws.Range("A6").Value = 112327043
$ws.Range("B6").Value = 96755
$ws.Range("E6").Value = 221952
$ws.Range("F6").Value = "Spindelblomster"
$ws.Range("G6").Value = "Neottia cordata"
$ws.Range("H6").Value = "(L.) Rich."
$ws.Range("Q6").Value = 641243
$ws.Range("R6").Value = 7163079

# Row 7
$ws.Range("A7").Value = 112327128
$ws.Range("B7").Value = 95701
$ws.Range("E7").Value = 221945
$ws.Range("F7").Value = "Revlummer"
$ws.Range("G7").Value = "Lycopodium annotinum"
$ws.Range("H7").Value = "L."
$ws.Range("Q7").Value = 641144
$ws.Range("R7").Value = 7163080

# Row 8
$ws.Range("A8").Value = 112327587
$ws.Range("B8").Value = 96652
$ws.Range("E8").Value = 219790
$ws.Range("F8").Value = "Fläcknycklar"
$ws.Range("G8").Value = "Dactylorhiza maculata"
$ws.Range("H8").Value = "(L.) Soó"
$ws.Range("Q8").Value = 641128

# Row 9
$ws.Range("A9").Value = 112327586
$ws.Range("B9").Value = 96652
$ws.Range("Q9").Value = 641118

# Row 10
$ws.Range("A10").Value = 112327584
$ws.Range("B10").Value = 96652
$ws.Range("Q10").Value = 641245

# Row 11
$ws.Range("A11").Value = 112327352
$ws.Range("B11").Value = 77650
$ws.Range("Q11").Value = 641114
